$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 281.74805
$ws.Range("B3").Value = 19.75
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2025-11-24"
$ws.Range("C3").ClearFormats()
$ws.Range("E3").Value = 0.0164
$ws.Range("B4").Value = 95.26000000000001
$ws.Range("B5").Value = 294.9145
$ws.Range("B6").Value = 1.56667
$ws.Range("B7").Value = 953.6774
$ws.Range("E7").Value = 0.0075
$ws.Range("B8").Value = 399.9678
$ws.Range("E8").Value = 0.0062
$ws.Range("B9").Value = 94.68273000000001
$ws.Range("E9").Value = 0.0171
$ws.Range("B10").Value = 2.9875
$ws.Range("B12").Value = 114.33333
$ws.Range("E12").Value = 0.0241
$ws.Range("B13").Value = 84.81476000000001
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "2026-01-02"
$ws.Range("C13").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "Feb 11, 2026"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = 0.0215
$ws.Range("B14").Value = 91.48
$ws.Range("E14").Value = 0.0341
$ws.Range("B15").Value = 3.625
$ws.Range("B16").Value = 132.50357
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "2025-12-15"
$ws.Range("C16").ClearFormats()
$ws.Range("E16").Value = 0.0147
$ws.Range("E17").Value = 0.0109
$ws.Range("B18").Value = 144.5
$ws.Range("B20").Value = 16.93418
$ws.Range("E20").Value = 0.0359
$ws.Range("B22").Value = 74.15385000000001
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "Jan 27, 2026"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = 0.008500000000000001
$ws.Range("B23").Value = 26.5
$ws.Range("E23").Value = 0.0247
$ws.Range("B24").Value = 27.7375
$ws.Range("E24").Value = 0.0475
$ws.Range("B25").Value = 463.54266
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "2025-11-27"
$ws.Range("C25").ClearFormats()
$ws.Range("E25").Value = 0.024
$ws.Range("E26").Value = 0.0218
$ws.Range("B28").Value = 7.233
$ws.Range("B29").Value = 37.31579
$ws.Range("B30").Value = 625.4096
$ws.Range("E30").Value = 0.0077
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "Jan 23, 2026"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = 0.0269
$ws.Range("B32").Value = 103.57
$ws.Range("E32").Value = 0.0116
$ws.Range("B33").Value = 134.65277
$ws.Range("B34").Value = 19.25054
$ws.Range("B35").Value = 47.125
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "2025-11-24"
$ws.Range("C35").ClearFormats()
$ws.Range("E35").Value = 0.014400001
$ws.Range("B37").Value = 6.6
$ws.Range("B39").Value = 49.6
$ws.Range("B40").Value = 47.23529
$ws.Range("E40").Value = 0.025
$ws.Range("B41").Value = 287.58334
$ws.Range("E41").Value = 0.0113
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "Jan 28, 2026"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = 0.0253
$ws.Range("B43").Value = 13.16667
$ws.Range("B44").Value = 392.9278
$ws.Range("B45").Value = 341.90668
$ws.Range("E45").Value = 0.012
$ws.Range("B47").Value = 7.2
$ws.Range("B48").Value = 399.9678
$ws.Range("E48").Value = 0.0062
$ws.Range("E49").Value = 0.0278
$ws.Range("B50").Value = 62.66667
$ws.Range("E50").Value = 0.050300002
$ws.Range("B51").Value = 12.125
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "2025-11-24"
$ws.Range("C51").ClearFormats()
$ws.Range("E51").Value = 0.0454
$ws.Range("B53").Value = 128.72
$ws.Range("E53").Value = 0.0355
$ws.Range("B56").Value = 16.30909
$ws.Range("E56").Value = 0.066199996
$ws.Range("B57").Value = 628.75
$ws.Range("E57").Value = 0.0252
$ws.Range("B58").Value = 12.515
$ws.Range("E58").Value = 0.0463
$ws.Range("B59").Value = 65.35278
$ws.Range("E59").Value = 0.0299
$ws.Range("B60").Value = 206
$ws.Range("B61").Value = 802.5263
$ws.Range("E61").Value = 0.0202
$ws.Range("B62").Value = 14.99999
$ws.Range("B63").Value = 34.23158
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "2025-12-26"
$ws.Range("C63").ClearFormats()
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "Jan 29, 2026"
$ws.Range("D63").ClearFormats()
$ws.Range("E63").Value = 0.0216
$ws.Range("B64").Value = 656.50555
$ws.Range("E64").Value = 0.0057
$ws.Range("B65").Value = 12.66667
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "Nov 6, 2025"
$ws.Range("D65").ClearFormats()
$ws.Range("B66").Value = 132.7711
$ws.Range("E66").Value = 0.031600002
$ws.Range("B67").Value = 15.47647
$ws.Range("E67").Value = 0.0538
$ws.Range("B68").Value = 50.41667
$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "2025-12-10"
$ws.Range("C68").ClearFormats()
$ws.Range("E68").Value = 0.0232
$ws.Range("B69").Value = 29.03958
$ws.Range("E69").Value = 0.0682
$ws.Range("B70").Value = 4.16538
$ws.Range("B71").Value = 42.04348
$ws.Range("E71").Value = 0.0094
$ws.Range("B72").Value = 84.00526000000001
$ws.Range("B73").Value = 233.42857
